$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the columns that were dropped from the import template -----
#   F - "Do you want to show the operator name on a buyer research ? (Yes or No)"
#   P - "Sale price amount"
#   Q - "Sale price currency"
# Delete right-to-left so the earlier column letters stay valid while we work.
$ws.Range("Q1").EntireColumn.Delete()
$ws.Range("P1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()

# --- Remove the trailing blank rows below the header table (old rows 11-54)
$ws.Range("A11:A54").EntireRow.Delete()

# --- Row heights were retightened along with the layout cleanup --------
$ws.Rows("1:2").RowHeight = 22.05
$ws.Rows("3:9").RowHeight = 12.8
$ws.Rows("10").RowHeight = 46.25

# --- Refresh the active selection / scroll position ---------------------
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F4").Select()

Write-Output ("Dimension after edit: " + $ws.UsedRange.Address())
